# TEXAS_2016.xlsx cleanup script
# 1. Rename header columns to short machine-friendly codes.
# 2. Title-case the Spanish connector words ("de", "del", "el", "la", "las",
#    "los", "y") that appear after the first word inside state/municipality
#    names (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga").
# 3. Drop the trailing metadata/footer rows (source notes) below the data,
#    shrinking the used range from A1:D1954 to A1:D1948.

function TitleCaseConnectors($s) {
    $connectors = @('de', 'del', 'el', 'la', 'las', 'los', 'y')
    $words = $s -split ' '
    for ($i = 1; $i -lt $words.Length; $i++) {
        if ($connectors -contains $words[$i]) {
            $w = $words[$i]
            $words[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1)
        }
    }
    return [string]::Join(' ', $words)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename headers (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Step 2: title-case connector words in columns A and B, rows 2-1948 ---
$lastRow = 1948
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null) {
        $cellA.Value = TitleCaseConnectors($valA)
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null) {
        $cellB.Value = TitleCaseConnectors($valB)
    }
}

# --- Step 3: delete trailing metadata rows 1950-1954 ---
$ws.Range("A1950:D1954").EntireRow.Delete()
